# Implementacion de disenos azul/amarillo version 1.0
# Applies the content edits for the new reporting period to the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (rows 6-11) ---------------------------------------------
$ws.Range("A7").Value = "GELOVER MANZO HUGO "
$ws.Range("E7").Value = "ING. ELECTROMECÁNICO"

$ws.Range("E9").Value = "01 DE SEPTIEMBRE DE 2019"
$ws.Range("H9").Value = "NO"

$ws.Range("A11").Value = "2024-1 (MARZO-AGOSTO 2024)"
$ws.Range("F11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0

# --- "Asignacion de Horas Frente a Grupo" table (rows 15-22) --------------
# Row 15: clear the career column, new subject row.
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = "CÁLCULO INTEGRAL (ACF-0902)"
$ws.Range("C15").Value = " "
$ws.Range("D15").Value = "5"
$ws.Range("E15").Value = "07:00 09:00"
$ws.Range("F15").Value = " "
$ws.Range("G15").Value = "07:00 09:00"
$ws.Range("H15").Value = " "
$ws.Range("I15").Value = "07:00 08:00"
$ws.Range("J15").Value = " "

# Row 16
$ws.Range("A16").Value = "INDUSTRIAL"
$ws.Range("B16").Value = "PROCESOS DE FABRICACIÓN (INC-1023)"
$ws.Range("C16").Value = "1401"
$ws.Range("D16").Value = "4"
$ws.Range("E16").Value = "09:00 11:00"
$ws.Range("G16").Value = "09:00 11:00"

# Row 17
$ws.Range("A17").Value = "INDUSTRIAL"
$ws.Range("B17").Value = "PROCESOS DE FABRICACIÓN (INC-1023)"
$ws.Range("C17").Value = "1451"
$ws.Range("D17").Value = "4"
$ws.Range("E17").Value = "15:00 17:00"
$ws.Range("G17").Value = "15:00 17:00"

# Row 18
$ws.Range("A18").Value = "INDUSTRIAL"
$ws.Range("B18").Value = "PROCESOS DE FABRICACIÓN (INC-1023)"
$ws.Range("C18").Value = "1481"
$ws.Range("D18").Value = "2"
$ws.Range("E18").Value = "11:00 13:00"

# Row 19
$ws.Range("A19").Value = "INDUSTRIAL"
$ws.Range("B19").Value = "PROCESOS DE FABRICACIÓN (INC-1023)"
$ws.Range("C19").Value = "1481"
$ws.Range("D19").Value = "2"
$ws.Range("E19").Value = "13:00 15:00 H.T. "

# Row 20
$ws.Range("A20").Value = "INDUSTRIAL"
$ws.Range("B20").Value = "ADMINISTRACIÓN DEL MANTENIMIENTO (INC-1004)"
$ws.Range("C20").Value = "1601"
$ws.Range("D20").Value = "4"
$ws.Range("G20").Value = "11:00 13:00"
$ws.Range("I20").Value = "11:00 13:00"

# Row 21
$ws.Range("A21").Value = "INDUSTRIAL"
$ws.Range("B21").Value = "ADMINISTRACIÓN DEL MANTENIMIENTO (INC-1004)"
$ws.Range("C21").Value = "1681"
$ws.Range("D21").Value = "3"
$ws.Range("G21").Value = "13:00 15:00 H.T. "
$ws.Range("I21").Value = "10:00 11:00"

# Row 22
$ws.Range("A22").Value = "INDUSTRIAL"
$ws.Range("B22").Value = "ADMINISTRACIÓN DEL MANTENIMIENTO (INC-1004)"
$ws.Range("C22").Value = "1681"
$ws.Range("D22").Value = "1"
$ws.Range("I22").Value = "13:00 14:00"

# Row 23 (TOTAL row for the table above)
$ws.Range("D23").Value = "25"

# --- "Asignacion de Horas de Descarga para otras Actividades" (rows 26-27) -
# Both activities were removed for this period.
$ws.Range("A26:D26").Value = ""
$ws.Range("E26").Value = " "
$ws.Range("A27:D27").Value = ""
$ws.Range("E27").Value = " "

# --- "Asignacion de Horas de Presidente y Secretario de Academia" total ---
$ws.Range("D34").Value = "0"

# --- Footer totals / application data --------------------------------------
$ws.Range("D39").Value = 0
$ws.Range("G40").Value = "01/03/2025"

# --- Signature block (rows 44-46) ------------------------------------------
$ws.Range("A44").Value = ""
$ws.Range("G44").Value = "GELOVER MANZO HUGO "

$ws.Range("A45").Value = ""

$ws.Range("A46").Value = "M. EN R.I. VIANCA LISSETH PEREZ CRUZ"
